$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows (11-13) for the "html" file-type code, mirroring the
# existing txt/xml/json pattern across the eng/ara/fra languages.

# Column A: code (same value "html" for all three new rows).
$ws.Range("A11").Value = "html"
$ws.Range("A12").Value = "html"
$ws.Range("A13").Value = "html"

# Column B: localized description - entered ara first, then eng, then fra.
$ws.Range("B12").Value = "ملف html"
$ws.Range("B11").Value = "html file"
$ws.Range("B13").Value = "Fichier html"

# Column C: lang_code
$ws.Range("C11").Value = "eng"
$ws.Range("C12").Value = "ara"
$ws.Range("C13").Value = "fra"

# Column D: is_active flag, left-aligned like the existing rows.
$ws.Range("D11").Value = $true
$ws.Range("D12").Value = $true
$ws.Range("D13").Value = $true
$ws.Range("D11:D13").HorizontalAlignment = -4131

# Column E: cr_by
$ws.Range("E11").Value = "superadmin"
$ws.Range("E12").Value = "superadmin"
$ws.Range("E13").Value = "superadmin"

# Column F: cr_dtimes
$ws.Range("F11").Value = "now()"
$ws.Range("F12").Value = "now()"
$ws.Range("F13").Value = "now()"

# Mirror the selection left behind in the source workbook after entering the
# data (selecting the remainder of the sheet starting at column G).
$ws.Range("G1:XFD1048576").Select()
